$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 15
$ws.Range("E15").Value = 149
$ws.Range("F15").Value = 77
$ws.Range("H15").Value = 77

# Row 36
$ws.Range("F36").Value = 38
$ws.Range("H36").Value = 38

# Row 38
$ws.Range("E38").Value = 65

# Row 46
$ws.Range("E46").Value = 25

# Row 47
$ws.Range("E47").Value = 53
$ws.Range("F47").Value = 32
$ws.Range("H47").Value = 32

# Row 61
$ws.Range("E61").Value = 27

# Row 63
$ws.Range("E63").Value = 29

# Row 77
$ws.Range("E77").Value = 50

# Row 81
$ws.Range("E81").Value = 15
$ws.Range("F81").Value = 6
$ws.Range("H81").Value = 6
